$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.229.95"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").Value = "1.583.59"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.90"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").Value = "0.498"
$ws.Range("E6").Value = "  -2.66%  "

$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "0.0611"
$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("E9").Value = "  -0.57%  "

$ws.Range("D10").Value = "19.53"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.585.92"
$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("E15").Value = "  -1.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").Value = "26.218.17"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").Value = "207.89"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  -1.13%  "

$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("E30").Value = "  -1.02%  "

$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("D34").Value = "1.278.70"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").Value = "0.608"
$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0166"
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "1.05"
$ws.Range("E39").Value = "  -10.08%  "

$ws.Range("D40").Value = "0.817"
$ws.Range("E40").Value = "  -1.51%  "

$ws.Range("E41").Value = "  +3.29%  "

$ws.Range("D42").Value = "0.768"
$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("E43").Value = "  -3.13%  "

$ws.Range("D44").Value = "62.31"

$ws.Range("D45").Value = "1.719.15"
$ws.Range("E45").Value = "  -1.14%  "

$ws.Range("D46").Value = "89.13"
$ws.Range("E46").Value = "  -1.68%  "

$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("E49").Value = "  -2.03%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.43"
$ws.Range("E51").Value = "  +0.17%  "
